$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0. Word re-drops the "_GoBack" last-edit-position bookmark at the
#    very start of the document content (this happens naturally
#    whenever the last editing activity before save is back at the
#    top of the document). Recreate the bookmark collapsed right at
#    the very start of the body text.
# ------------------------------------------------------------------
$docStart = $d.Range(0, 0)
$docStart.InsertBefore("z")
$tmpRng = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $tmpRng) | Out-Null
$d.Range(0, 1).Delete() | Out-Null

# ------------------------------------------------------------------
# 1. Fix wording: "formando por" -> "formado por"
# ------------------------------------------------------------------
$d.Content.Find.Execute("formando por", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "formado por", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Italicize the genus names Penicillium, Aspergillus and Rhizopus
#    in "...pertenecen a los géneros Penicillium, Aspergillus y Rhizopus."
#    ("Penicillium" is unique in the doc; "Aspergillus" and "Rhizopus"
#    each occur twice, so we anchor on the unique surrounding phrase
#    and then shrink the found range down to just the genus word.)
# ------------------------------------------------------------------

# -- Penicillium (already unique) --
$rng = $d.Content
$rng.Find.Execute("Penicillium", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Font.Italic = $true

# -- Aspergillus (the one immediately followed by " y ") --
$rng = $d.Content
$rng.Find.Execute("Aspergillus y ", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.MoveEnd(1, -3) | Out-Null   # drop the trailing " y " (3 chars) -> just "Aspergillus"
$rng.Font.Italic = $true

# -- Rhizopus (the one immediately preceded by "Aspergillus y ") --
$rng = $d.Content
$rng.Find.Execute("Aspergillus y Rhizopus", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.MoveStart(1, 14) | Out-Null # drop the leading "Aspergillus y " (14 chars) -> just "Rhizopus"
$rng.Font.Italic = $true
